$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3-4: capitalize "center"/"봉사자" -> "Center"/"Volunteer" in the use-case name column ---
$ws.Range("B3").Value = "Center 회원가입"
$ws.Range("B4").Value = "Volunteer 회원가입"

# --- Rows 5-7: a brand new "로그인" (Login) use-case is inserted, shifting/renaming the
#     remaining rows. Re-write IDs, names and descriptions for rows 5-11. ---
$ws.Range("A5").Value = "UC-04"
$ws.Range("B5").Value = "로그인"
$ws.Range("C5").Value = "center가 로그인 하는 기능을 설정"

$ws.Range("A6").Value = "UC-05"
$ws.Range("B6").Value = "Center 로그인"
$ws.Range("C6").Value = "Volunteer가 로그인 하는 기능을 설정"

$ws.Range("A7").Value = "UC-06"
$ws.Range("B7").Value = "Volunteer 로그인"
$ws.Range("C7").Value = "센터 로그인, 봉사자 로그인 선택"

$ws.Range("A8").Value = "UC-07"
$ws.Range("B8").Value = "Volunteer 구인"
$ws.Range("C8").ClearContents()

$ws.Range("A9").Value = "UC-08"
$ws.Range("B9").Value = "Volunteer 예약"
$ws.Range("C9").ClearContents()

$ws.Range("A10").Value = "UC-09"
$ws.Range("B10").Value = "Center 정보 업데이트"
$ws.Range("C10").ClearContents()

$ws.Range("A11").Value = "UC-10"
$ws.Range("B11").Value = "Volunteer 회원정보 변경"
$ws.Range("C11").ClearContents()

# --- Column B is a bit wider now ---
$ws.Columns.Item(2).ColumnWidth = 22.36

# --- Cursor/selection moved to C8 ---
$ws.Range("C8").Select()
